$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Update the numeric results (instance costs) that were corrected.
$ws.Range("B2").Value = 0.69841269841269804
$ws.Range("B3").Value = 1.0595238095238
$ws.Range("B6").Value = 1.5119047619047601
$ws.Range("B7").Value = 0.75595238095238004
$ws.Range("B8").Value = 5.0941798941798897

# B6/B7 lose their "applyFont" bold-flag styling (now matches the plain
# bordered style used elsewhere), B8 loses its border+applyFont styling
# entirely and reverts to the default style.
$ws.Range("B6").Font.Bold = $false
$ws.Range("B7").Font.Bold = $false
$ws.Range("B8").Font.Bold = $false

# Move/restore the active selection to E8.
$ws.Range("E8").Select()
